$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace each cell's previous text with the new value. Using Replace()
# scoped to a single-cell Range (rather than assigning .Value) makes the
# engine recompute the shared-string table, pruning entries that become
# unused - matching how Excel rewrites xl/sharedStrings.xml on save.

$ws.Range("A2").Replace("62b68573d5bacf79bd396d8d", "62b74c1db2ceda5fa906598f")
$ws.Range("B2").Replace("Week 1", "Thomas")

$ws.Range("A3").Replace("62b6857121e723731dcb07d6", "62b74c1ebfc68954fbf6f5e5")
$ws.Range("B3").Replace("Week 1", "Ying")

$ws.Range("A4").Replace("62b6856f1b0c9b18e5e84416", "62b74c20440fb27b11724a3d")
$ws.Range("B4").Replace("Week 1", "Esther")

$ws.Range("A5").Replace("62b6856de3108247dbcb359b", "62b74c21e17fdb80e8513e7a")
$ws.Range("B5").Replace("Week 1", "Zaur")

$ws.Range("A6").Replace("62b666a6c9e1e9804b1972b4", "62b74c23e2197787f1b7e3b4")
$ws.Range("B6").Replace("To Do", "Diahandra")

# Rows 7-8 no longer carry a person/card id in column A, and column B is
# blanked out (the row is kept, just emptied).
$ws.Range("A7").Clear()
$ws.Range("B7").Replace("Doing", "")

$ws.Range("A8").Clear()
$ws.Range("B8").Replace("Done", "")
